$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 64, shifting existing rows 64-72 down to 65-73
$ws.Rows("64:64").Insert()

# Populate the new row 64 with this week's data (weekly update for
# Terminal Hortofrutícola Agro Chillán - Haba)
$ws.Range("A64").Value = 7
$ws.Range("B64").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C64").Value = "Ñuble"
$ws.Range("D64").Value = 44918
$ws.Range("E64").Value = 16
$ws.Range("F64").Value = 100112026
$ws.Range("G64").Value = "Haba"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 100
$ws.Range("K64").Value = 15000
$ws.Range("L64").Value = 16000
$ws.Range("M64").Value = 15500
$ws.Range("N64").Value = "$/saco 25 kilos"
$ws.Range("O64").Value = "Provincia de Diguillín"
$ws.Range("P64").Value = 620
$ws.Range("Q64").Value = 25
$ws.Range("R64").Value = "Hortaliza"
